# "Next Steps" slide: rework bullet text adjustments.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

# Paragraph 1: "Further rework the draft (structure and application examples)"
#           -> "Further rework the draft (content, structure and examples)"
# Two-step set avoids the host's word-level diff from leaving stale runs
# around the unchanged words, keeping the paragraph a single run.
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "."
$para1 = $tr.Paragraphs(1, 1)
$para1.Text = "Further rework the draft (content, structure and examples)"

# Paragraph 2: was 3 runs ("Investigate 2" / "nd" superscript / " JWS serialization option ")
#           -> single run "Alignment in BRSKI design team calls"
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "."
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Alignment in BRSKI design team calls"
